$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  46000, 11199.4657474991, 10427.0727614257, 11232.26, 6860.98814747214, 252.325037870744),
    @(3,  46001, 11451.4431360528, 10497.4668600735, 11232.26, 7148.51755574849, 267.238517325915),
    @(4,  46002, 11485.7393482061, 10490.7778761424, 11232.26, 7173.44431758768, 267.998424738752),
    @(5,  46003, 11382.9199775364, 9679.49501174235,  11232.26, 7112.70642018032, 231.664226330111),
    @(6,  46004, 4034.69176065649, 6526.43696927929,  11232.26, 6392.20995135107, 70.2661216929316),
    @(7,  46005, 8386.73995792576, 8957.64025669508,  11232.26, 7357.50840784026, 211.787027688972),
    @(8,  46006, 9746.03784374469, 9889.58170276042,  11232.26, 8042.21065652037, 279.1471816367),
    @(9,  46007, 9746.03784374469, 9991.1487090866,   11232.26, 8042.21065652037, 283.379140233624),
    @(10, 46008, 9746.03784374469, 9744.82225048386,  11232.26, 8042.21065652037, 273.115537791843),
    @(11, 46009, 9746.03784374469, 9719.90506739582,  11232.26, 8042.21065652037, 272.077321829841),
    @(12, 46010, 9746.03784374469, 8955.30978858843,  11232.26, 8042.21065652037, 240.219185212867),
    @(13, 46011, 8560.57524882407, 8794.01074684497,  11232.26, 7650.27361136198, 217.167681591956),
    @(14, 46012, 8469.04674334209, 8870.38913361643,  11232.26, 7641.87559525584, 220.000197036345),
    @(15, 46013, 9566.18303866457, 9512.59216458739,  11232.26, 8113.08526256509, 266.39239279802)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
